{"js": "// Find the existing \"Implemented Services using Elixir 1.3\" bullet paragraph\n// and insert a new bullet paragraph right after it, reading\n// \"Implemented Services using Clojure 1.8\" \u2014 matching the same list\n// (numId 2) and run formatting (Arial, 10pt).\nconst body = context.document.body;\nconst searchResults = body.search(\"Implemented Services using Elixir 1.3\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find target paragraph text to anchor the insert.\");\n}\n\nconst anchorParagraph = searchResults.items[0].paragraphs.getFirst();\nconst newParagraph = anchorParagraph.insertParagraph(\n  \"Implemented Services using Clojure 1.8\",\n  \"After\"\n);\n\n// Match the formatting of the existing achievement bullets: Arial, 10pt (sz 20 half-points).\nnewParagraph.font.set({ name: \"Arial\", size: 10 });\n\nawait context.sync();\n", "ps1": "# Insert a new \"Implemented Services using Clojure 1.8\" bullet right after\n# the existing \"Implemented Services using Elixir 1.3\" bullet in the\n# Achievements list, matching that bullet's list/indent/font formatting.\n\n$d = $word.ActiveDocument\n\n$targetText = \"Implemented Services using Elixir 1.3\"\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\n\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  $pText = $p.Range.Text.TrimEnd(\"`r\", \"`a\")\n  if ($pText -eq $targetText) {\n    $anchorIndex = $i\n    break\n  }\n}\n\nif ($anchorIndex -eq -1) {\n  throw \"Could not find anchor paragraph '$targetText'\"\n}\n\n$anchorPara = $d.Paragraphs.Item($anchorIndex)\n\n# InsertParagraphAfter() clones the anchor paragraph's formatting (pPr/rPr,\n# including the bullet's numId/ilvl and the Arial/10pt run font) onto a new,\n# empty paragraph placed right after it.\n$anchorPara.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item($anchorIndex + 1)\n$newPara.Range.Text = \"Implemented Services using Clojure 1.8\"\n"}
